# Adding "Average Number of Words Per Sentence" summary row to the
# Categorized Records worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label + formula in row 54 (two blank rows below the last data row, 51)
$ws.Range("B54").Value = "Average Number of Words Per Sentence"
$ws.Range("C54").Formula = "=SUM(C2:C51)/50"

# Give the new summary cells the same thin box border used elsewhere in the
# sheet (e.g. the header row), without the bold font/alignment formatting.
$rng = $ws.Range("B54:C54")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Resize column B to fit the new (longer) label text.
$ws.Columns(2).AutoFit()

# Update the on-screen selection / scroll position to match where the user
# ended up after adding the new row.
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$null = $ws.Range("B59").Select()
